$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as plain text without altering its style,
# by temporarily marking the cell as Text format, assigning the value,
# then resetting the style back to Normal (removes the style index again).
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "69.288.71"
$ws.Range("E2").Value = "  -2.91%  "
Set-TextValue "D3" "3.496.98"
$ws.Range("E3").Value = "  -2.43%  "
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  -0.10%  "
Set-TextValue "D5" "607.73"
$ws.Range("E5").Value = "  +3.95%  "
Set-TextValue "D6" "185.25"
$ws.Range("E6").Value = "  -1.28%  "
Set-TextValue "D7" "0.625"
$ws.Range("E7").Value = "  +0.28%  "
Set-TextValue "D8" "1.00"
$ws.Range("E8").Value = "  +0.02%  "
Set-TextValue "D9" "0.209"
$ws.Range("E9").Value = "  -3.33%  "
Set-TextValue "D10" "0.652"
$ws.Range("E10").Value = "  -0.59%  "
Set-TextValue "D11" "53.19"
$ws.Range("E11").Value = "  -2.90%  "
$ws.Range("E12").Value = "  -4.34%  "
Set-TextValue "D13" "9.61"
$ws.Range("E13").Value = "  +0.18%  "
Set-TextValue "D14" "4.061.66"
Set-TextValue "D15" "612.17"
$ws.Range("E15").Value = "  +8.02%  "
$ws.Range("E16").Value = "  +2.27%  "
$ws.Range("E17").Value = "  -1.83%  "
Set-TextValue "D18" "69.387.48"
$ws.Range("E18").Value = "  -2.71%  "
$ws.Range("E19").Value = "  -2.34%  "
$ws.Range("E20").Value = "  -0.32%  "
$ws.Range("E21").Value = "  -2.66%  "
Set-TextValue "D22" "17.53"
$ws.Range("E22").Value = "  -0.28%  "
Set-TextValue "D23" "104.71"
$ws.Range("E23").Value = "  +10.36%  "
$ws.Range("E24").Value = "  +1.34%  "
Set-TextValue "D25" "5.01"
$ws.Range("E25").Value = "  -0.98%  "
$ws.Range("E26").Value = "  +1.53%  "
Set-TextValue "D27" "10.87"
$ws.Range("E27").Value = "  -4.10%  "
Set-TextValue "D28" "9.87"
$ws.Range("E28").Value = "  +7.02%  "
Set-TextValue "D29" "33.65"
$ws.Range("E29").Value = "  +3.02%  "
Set-TextValue "D30" "6.98"
$ws.Range("E30").Value = "  -4.80%  "
$ws.Range("E31").Value = "  +0.48%  "
Set-TextValue "D33" "63.39"
$ws.Range("E33").Value = "  -1.29%  "
Set-TextValue "D34" "3.73"
$ws.Range("E34").Value = "  +13.61%  "
$ws.Range("B35").Value = "Dai"
$ws.Range("C35").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D35" "1.00"
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D36" "3.11"
$ws.Range("E36").Value = "  -9.07%  "
Set-TextValue "D37" "524.51"
$ws.Range("E37").Value = "  -5.54%  "
Set-TextValue "D38" "0.395"
$ws.Range("E38").Value = "  -6.68%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D39" "3.57"
$ws.Range("E39").Value = "  +3.36%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D40" "3.545.09"
$ws.Range("E40").Value = "  +0.34%  "
Set-TextValue "D41" "36.53"
$ws.Range("E41").Value = "  -3.49%  "
$ws.Range("E42").Value = "  +2.81%  "
$ws.Range("E43").Value = "  -5.25%  "
Set-TextValue "D44" "0.0457"
$ws.Range("E44").Value = "  +1.27%  "
Set-TextValue "D45" "2.95"
$ws.Range("E45").Value = "  -0.50%  "
$ws.Range("E46").Value = "  +3.78%  "
$ws.Range("E47").Value = "  -3.99%  "
Set-TextValue "D48" "8.86"
$ws.Range("E48").Value = "  -5.84%  "
$ws.Range("E49").Value = "  +0.24%  "
Set-TextValue "D50" "132.53"
$ws.Range("E50").Value = "  -2.53%  "
$ws.Range("E51").Value = "  -9.48%  "
